# Add "Unsuitable exception type" check as a new row (row 32) on the
# "Workflow" sheet, following the same layout/style as the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

$newRow = 32

# Copy formatting from the row above (row 31, columns A:G only) so the new
# row matches the existing table's look (borders, fill, wrap, etc.) without
# touching the rest of the (unused) row.
$ws.Range("A31:G31").Copy()
$ws.Range("A32:G32").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Column layout: A=Run, B=Issue, C=Check Filename, D=Arguments, E=Action,
# F=Explanation, G=Suggestion
$ws.Cells.Item($newRow, 1).Value = "No"
$ws.Cells.Item($newRow, 2).Value = "Unsuitable exception type"
$ws.Cells.Item($newRow, 3).Value = "Checks\Custom\UnsuitableExceptionType.xaml"
# D32 (Arguments) intentionally left blank for this check.
$ws.Cells.Item($newRow, 5).Value = "Fix"
$ws.Cells.Item($newRow, 6).Value = "When throwing exceptions, it is recommended to properly distinguish between application-originated and business-originated exceptions. The type of the exception to be thrown or caught should be as specific as possible, and Exception and ApplicationException should be avoided."
$ws.Cells.Item($newRow, 7).Value = "Use specific exception types and avoid using generic types such as Exception and ApplicationException."

$ws.Rows.Item($newRow).RowHeight = 87

# Extend the existing data validation ranges so the new row is included,
# matching the updated sqref values A2:A32 and E11:E32.
$ws.Range("E11:E31").Validation.Delete()
$ws.Range("E11:E32").Validation.Add(3, 1, 1, '"Fix, Double check"')

$ws.Range("A2:A31").Validation.Delete()
$ws.Range("A2:A32").Validation.Add(3, 1, 1, '"Yes, No"')
